# Updates cryptos list: refreshed Price/Volume(1h) figures, plus a
# rank swap between PaxDollar and NEARProtocol (rows 49/50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.311.03"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "1.820.60"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'313.63"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.4646"
$ws.Range("E7").Value = "  +4.63%  "
$ws.Range("D8").Value = "'0.3773"
$ws.Range("E8").Value = "  +2.62%  "
$ws.Range("D9").Value = "'0.07428"
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("D10").Value = "'0.8710"
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("D11").Value = "'20.65"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "1.825.49"
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").Value = "'6.684"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").Value = "'5.406"
$ws.Range("D15").Value = "'0.07103"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "'92.15"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "'0.000008769"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'14.94"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("D21").Value = "27.313.75"
$ws.Range("E21").Value = "  +2.35%  "
$ws.Range("D22").Value = "'5.315"
$ws.Range("E22").Value = "  +3.20%  "
$ws.Range("D23").Value = "'10.95"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").Value = "2.049.90"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("D25").Value = "'1.936"
$ws.Range("E25").Value = "  -2.26%  "
$ws.Range("D26").Value = "'151.70"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "'2.249"
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("D28").Value = "'18.60"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").Value = "'5.297"
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("D30").Value = "'117.26"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "'0.08924"
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("D32").Value = "'0.7813"
$ws.Range("E32").Value = "  +6.04%  "
$ws.Range("D33").Value = "'1.184"
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("D34").Value = "'4.522"
$ws.Range("E34").Value = "  +2.02%  "
$ws.Range("D35").Value = "'2.920"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").Value = "'1.001"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'1.098"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").Value = "'0.01971"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").Value = "'0.05247"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D40").Value = "'7.310"
$ws.Range("E40").Value = "  +5.02%  "
$ws.Range("D41").Value = "'2.368"
$ws.Range("E41").Value = "  +20.54%  "
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").Value = "'2.888"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("D44").Value = "'0.1691"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("E45").Value = "  +2.17%  "
$ws.Range("D46").Value = "'0.5044"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "'10.49"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("D48").Value = "'105.58"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.673"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'1.000"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").Value = "'0.06328"
$ws.Range("E51").Value = "  +0.64%  "
